$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "TestCasesFlag" (sheet1): add TC_002..TC_005 rows
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("TestCasesFlag")

$ws1.Range("A3").Value = "TC_002_loginandVerify"
$ws1.Range("B3").Value = $false

$ws1.Range("A4").Value = "TC_003_loginandVerify"
$ws1.Range("B4").Value = $false

$ws1.Range("A5").Value = "TC_004_loginandVerify"
$ws1.Range("B5").Value = $false

$ws1.Range("A6").Value = "TC_005_loginandVerify"
$ws1.Range("B6").Value = $false

$ws1.Range("B2").Select()

# ---------------------------------------------------------------
# Sheet "Data" (sheet2): rework columns / contents
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Data")

# Column widths (approximate closest achievable values for
# width="23.36328125" / width="14.453125")
$ws2.Columns.Item(1).ColumnWidth = 22.5
$ws2.Columns.Item(2).ColumnWidth = 13.666666666666666

# Row 1 - headers
$ws2.Range("A1").Value = "TestCaseName"
$ws2.Range("B1").Value = "userid"
$ws2.Range("C1").Value = "password"
$ws2.Range("D1").Value = "FirstName"
$ws2.Range("E1").Value = "LastName"
$ws2.Range("F1").Value = "MidName "
$ws2.Range("G1").Value = "Add"
$ws2.Range("H1").Value = "Add2"

# Row 2
$ws2.Range("A2").Value = "TC_001_loginandVerify"
$ws2.Range("B2").Value = "admin"
$ws2.Range("C2").Value = "Password01"
$ws2.Range("D2").Value = "AML admin"
$ws2.Range("E2").Value = "UB"
$ws2.Range("F2").Value = "PRAKASH"
$ws2.Range("G2").Value = "OSLO"
$ws2.Range("H2").Value = "Done"

# Row 3
$ws2.Range("A3").Value = "TC_001_loginandVerify"
$ws2.Range("B3").Value = "Controller"
$ws2.Range("C3").Value = "Password02"
$ws2.Range("D3").Value = "AML controller"
$ws2.Range("E3").Value = "PRABIN"
$ws2.Range("F3").Value = "KAUR"
$ws2.Range("G3").Value = "EDINBURG"
$ws2.Range("H3").Value = "London"

# Row 4 - new row
$ws2.Range("A4").Value = "TC_003_loginandVerify"
$ws2.Range("B4").Value = "three"
$ws2.Range("C4").Value = "three"
$ws2.Range("D4").Value = "three"
$ws2.Range("E4").Value = "three"
$ws2.Range("F4").Value = "three"
$ws2.Range("G4").Value = "three"
$ws2.Range("H4").Value = "three"

$ws2.Activate()
$ws2.Range("C4").Select()
